$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.109.94"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "1.789.45"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("D4").Value = "'0.999"

$ws.Range("D5").Value = "'222.27"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'32.33"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "'0.0715"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "2.043.91"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "1.795.41"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").Value = "'10.92"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "'0.628"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").Value = "34.082.05"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").Value = "'4.18"
$ws.Range("E17").Value = "  -3.16%  "

$ws.Range("D18").Value = "'68.12"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").Value = "'244.31"
$ws.Range("E19").Value = "  -3.93%  "

$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "'10.83"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = "  -3.53%  "

$ws.Range("E24").Value = "  -1.18%  "

$ws.Range("D25").Value = "'159.08"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").Value = "'16.37"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "'7.07"
$ws.Range("E27").Value = "  -0.96%  "

$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -2.83%  "

$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").Value = "'3.68"
$ws.Range("E32").Value = "  -3.49%  "

$ws.Range("D33").Value = "'3.50"
$ws.Range("E33").Value = "  -3.29%  "

$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("D35").Value = "1.399.93"
$ws.Range("E35").Value = "  -2.89%  "

$ws.Range("D36").Value = "'0.654"
$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("E38").Value = "  -3.45%  "

$ws.Range("D39").Value = "'79.80"
$ws.Range("E39").Value = "  -6.19%  "

$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").Value = "'0.922"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D42").Value = "'2.72"
$ws.Range("E42").Value = "  -2.64%  "

$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.91"
$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "'0.0495"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "'107.67"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D48").Value = "1.944.60"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").Value = "'12.08"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("E51").Value = "  +0.85%  "
